$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title / date-range text (edit in place to preserve rich-text run styling) ---
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "11"

$dateCell = $ws.Range("C9")
$dateCell.Characters(46, 8).Text = "3/16/2025"
$dateCell.Characters(27, 8).Text = "3/10/2025"

# --- Update table values (rows 14-30) ---
$ws.Range("L14").Value = -75
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = -100
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = 0
$ws.Range("N15").Value = -88.235294117647
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 15.384615384615
$ws.Range("I16").Value = 48
$ws.Range("J16").Value = 39
$ws.Range("K16").Value = 23.076923076923
$ws.Range("L16").Value = 41.176470588235
$ws.Range("M16").Value = -4
$ws.Range("N16").Value = -69.032258064516
$ws.Range("C17").Value = 12
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = 0
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = 31.034482758620
$ws.Range("I17").Value = 88
$ws.Range("J17").Value = 79
$ws.Range("K17").Value = 11.392405063291
$ws.Range("L17").Value = 25.714285714285
$ws.Range("M17").Value = 193.333333333333
$ws.Range("N17").Value = -19.266055045871
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -71.428571428571
$ws.Range("I18").Value = 18
$ws.Range("J18").Value = 38
$ws.Range("K18").Value = -52.631578947368
$ws.Range("L18").Value = -30.769230769230
$ws.Range("M18").Value = 28.571428571428
$ws.Range("N18").Value = -78.823529411764
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 71.428571428571
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 21.428571428571
$ws.Range("I19").Value = 83
$ws.Range("J19").Value = 76
$ws.Range("K19").Value = 9.210526315789
$ws.Range("L19").Value = -9.782608695652
$ws.Range("M19").Value = 80.434782608695
$ws.Range("N19").Value = -36.641221374045
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 8
$ws.Range("J20").Value = 11
$ws.Range("K20").Value = -27.272727272727
$ws.Range("L20").Value = -66.666666666666
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -89.610389610389
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = 6.896551724137
$ws.Range("F21").Value = 98
$ws.Range("G21").Value = 97
$ws.Range("H21").Value = 1.030927835051
$ws.Range("I21").Value = 248
$ws.Range("J21").Value = 246
$ws.Range("K21").Value = 0.813008130081
$ws.Range("L21").Value = -2.362204724409
$ws.Range("M21").Value = 68.707482993197
$ws.Range("N21").Value = -57.241379310344
$ws.Range("C23").Value = 12
$ws.Range("D23").Value = 8
$ws.Range("F23").Value = 30
$ws.Range("H23").Value = 20
$ws.Range("I23").Value = 73
$ws.Range("J23").Value = 71
$ws.Range("K23").Value = 2.816901408450
$ws.Range("L23").Value = 1.388888888888
$ws.Range("M23").Value = 87.179487179487
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 18
$ws.Range("F24").Value = 83
$ws.Range("G24").Value = 63
$ws.Range("H24").Value = 31.746031746031
$ws.Range("I24").Value = 210
$ws.Range("J24").Value = 157
$ws.Range("K24").Value = 33.757961783439
$ws.Range("L24").Value = 31.25
$ws.Range("M24").Value = 76.470588235294
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 300
$ws.Range("F25").Value = 34
$ws.Range("H25").Value = 142.857142857143
$ws.Range("I25").Value = 73
$ws.Range("J25").Value = 27
$ws.Range("K25").Value = 170.37037037037
$ws.Range("L25").Value = 114.705882352941
$ws.Range("C26").Value = 17
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = 21.428571428571
$ws.Range("I26").Value = 108
$ws.Range("J26").Value = 132
$ws.Range("K26").Value = -18.181818181818
$ws.Range("L26").Value = 1.886792452830
$ws.Range("M26").Value = -8.474576271186
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("J27").Value = 5
$ws.Range("K27").Value = -40
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -20
$ws.Range("I28").Value = 10
$ws.Range("J28").Value = 11
$ws.Range("K28").Value = -9.090909090909
$ws.Range("L28").Value = 0
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 2
$ws.Range("I29").Value = 2
$ws.Range("K29").Value = 100
$ws.Range("L29").Value = -60
$ws.Range("M29").Value = -60
$ws.Range("N29").Value = -92.307692307692
$ws.Range("C30").Value = 1
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 2
$ws.Range("K30").Value = 100
$ws.Range("L30").Value = -60
$ws.Range("M30").Value = -60
$ws.Range("N30").Value = -92

# --- Cells whose number format switches from General (placeholder text) to a real numeric format ---
$intFormatCells = @("D15", "G15", "D20", "C27", "F27", "C29", "C30")
foreach ($r in $intFormatCells) {
    $ws.Range($r).NumberFormat = "#,##0"
}

$decFormatCells = @("E15", "H15", "E20")
foreach ($r in $decFormatCells) {
    $ws.Range($r).NumberFormat = "#,##0.0;""-""#,##0.0"
}
